# Locate the paragraph containing the target phrase "unsorted heap of clothes"
# and rewrite it so that the trailing word "clothes" is spell-check wrapped
# and immediately followed by a new run of stray text "www.gmai".
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*heap of clothes*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Output "Target paragraph not found"
} else {
    $wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

    $xml = "<w:p $wns w:rsidR=`"005C6D3F`" w:rsidRPr=`"002A4D98`" w:rsidRDefault=`"005C6D3F`" w:rsidP=`"005C6D3F`">" +
        "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/></w:rPr></w:pPr>" +
        "<w:proofErr w:type=`"gramStart`"/>" +
        "<w:r w:rsidRPr=`"002A4D98`"><w:rPr><w:sz w:val=`"22`"/></w:rPr><w:t>unsorted</w:t></w:r>" +
        "<w:proofErr w:type=`"gramEnd`"/>" +
        "<w:r w:rsidRPr=`"002A4D98`"><w:rPr><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`"> heap of </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:rPr><w:sz w:val=`"22`"/></w:rPr><w:t>clothes</w:t></w:r>" +
        "<w:r><w:rPr><w:sz w:val=`"22`"/></w:rPr><w:t>www.gmai</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "</w:p>"

    $target.Range.InsertXML($xml) | Out-Null
    Write-Output ("Updated paragraph text: " + $target.Range.Text)
}
